$d = $word.ActiveDocument

# The page footer chrome scraped off the site - an empty paragraph,
# "Ver no Jupiter Salvar em pdf Salvar em docx" and the "(c) 2020 ..."
# copyright line - is being removed. Find the two text paragraphs, then
# widen the deletion to also swallow the blank paragraph immediately
# preceding the "Ver no Jupiter" line, and delete the whole span
# (through the end of the copyright paragraph, paragraph mark included)
# in one go.

$count = $d.Paragraphs.Count
$jupiterIndex = 0
$copyrightIndex = 0

for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $jupiterIndex = $i
    }
    if ($t -like "*Contact: luizeleno@usp.br*") {
        $copyrightIndex = $i
    }
}

$startPara = $d.Paragraphs.Item($jupiterIndex - 1)
$endPara = $d.Paragraphs.Item($copyrightIndex)

$r = $d.Range($startPara.Range.Start, $endPara.Range.End)
$r.Delete()
